{"js": "// Map of old math-fact strings to their new replacements, exactly as they\n// appear as standalone text runs inside the table cells of the document.\nconst replacements = [\n  [\"609\u00d72=1218\", \"399\u00d72=798\"],\n  [\"850\u00d76=5100\", \"888\u00d73=2664\"],\n  [\"117\u00d79=1053\", \"667\u00d76=4002\"],\n  [\"411\u00d77=2877\", \"957\u00d75=4785\"],\n  [\"490\u00d75=2450\", \"775\u00d76=4650\"],\n  [\"546\u00d79=4914\", \"555\u00d72=1110\"],\n  [\"460\u00d76=2760\", \"444\u00d79=3996\"],\n  [\"868\u00d75=4340\", \"162\u00d74=648\"],\n  [\"638\u00d75=3190\", \"972\u00d72=1944\"],\n  [\"520\u00d73=1560\", \"469\u00d79=4221\"],\n  [\"417\u00d76=2502\", \"312\u00d79=2808\"],\n  [\"625\u00d72=1250\", \"870\u00d74=3480\"],\n  [\"683\u00d79=6147\", \"816\u00d72=1632\"],\n  [\"612\u00d75=3060\", \"646\u00d75=3230\"],\n  [\"883\u00d78=7064\", \"156\u00d77=1092\"],\n  [\"248\u00d75=1240\", \"264\u00d79=2376\"],\n  [\"403\u00d78=3224\", \"786\u00d73=2358\"],\n  [\"776\u00d77=5432\", \"592\u00d78=4736\"],\n  [\"728\u00d79=6552\", \"781\u00d77=5467\"],\n  [\"642\u00d72=1284\", \"542\u00d79=4878\"],\n  [\"554\u00d79=4986\", \"809\u00d78=6472\"],\n  [\"192\u00d75=960\", \"905\u00d72=1810\"],\n  [\"412\u00d75=2060\", \"268\u00d75=1340\"],\n  [\"283\u00d74=1132\", \"761\u00d78=6088\"],\n  [\"265\u00d75=1325\", \"565\u00d77=3955\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each three-digit-times-one-digit math fact with its new value.\n# Every old string occurs exactly once in the document, inside a table cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"609\u00d72=1218\"; New = \"399\u00d72=798\" },\n    @{ Old = \"850\u00d76=5100\"; New = \"888\u00d73=2664\" },\n    @{ Old = \"117\u00d79=1053\"; New = \"667\u00d76=4002\" },\n    @{ Old = \"411\u00d77=2877\"; New = \"957\u00d75=4785\" },\n    @{ Old = \"490\u00d75=2450\"; New = \"775\u00d76=4650\" },\n    @{ Old = \"546\u00d79=4914\"; New = \"555\u00d72=1110\" },\n    @{ Old = \"460\u00d76=2760\"; New = \"444\u00d79=3996\" },\n    @{ Old = \"868\u00d75=4340\"; New = \"162\u00d74=648\" },\n    @{ Old = \"638\u00d75=3190\"; New = \"972\u00d72=1944\" },\n    @{ Old = \"520\u00d73=1560\"; New = \"469\u00d79=4221\" },\n    @{ Old = \"417\u00d76=2502\"; New = \"312\u00d79=2808\" },\n    @{ Old = \"625\u00d72=1250\"; New = \"870\u00d74=3480\" },\n    @{ Old = \"683\u00d79=6147\"; New = \"816\u00d72=1632\" },\n    @{ Old = \"612\u00d75=3060\"; New = \"646\u00d75=3230\" },\n    @{ Old = \"883\u00d78=7064\"; New = \"156\u00d77=1092\" },\n    @{ Old = \"248\u00d75=1240\"; New = \"264\u00d79=2376\" },\n    @{ Old = \"403\u00d78=3224\"; New = \"786\u00d73=2358\" },\n    @{ Old = \"776\u00d77=5432\"; New = \"592\u00d78=4736\" },\n    @{ Old = \"728\u00d79=6552\"; New = \"781\u00d77=5467\" },\n    @{ Old = \"642\u00d72=1284\"; New = \"542\u00d79=4878\" },\n    @{ Old = \"554\u00d79=4986\"; New = \"809\u00d78=6472\" },\n    @{ Old = \"192\u00d75=960\";  New = \"905\u00d72=1810\" },\n    @{ Old = \"412\u00d75=2060\"; New = \"268\u00d75=1340\" },\n    @{ Old = \"283\u00d74=1132\"; New = \"761\u00d78=6088\" },\n    @{ Old = \"265\u00d75=1325\"; New = \"565\u00d77=3955\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
